$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fix the model name typo: "DT1" -> "DT"
$ws.Range("B5").Value = "Voting_LogReg_Adaboost_w_DT_Xgboost"

# Update the active selection on the sheet
$ws.Range("C8").Select()
